$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows where the "Approved/Rejected" column (I) was "Rejected" with a
# "ReasonToReject" (J) of "Nil" -> change to "Approved" and clear the reason.
$rows = @(13, 14, 17, 18, 19, 20, 21, 22)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "Approved"
    $ws.Cells.Item($r, 10).ClearContents()
}

# Update the view: move the active selection to A22 (this also resets the
# previously scrolled top-left cell back to the sheet's default).
$ws.Range("A22").Select()
